# Regenerate merged AHB files
#
# 1. Rename the "_old"/"_new" suffixed column headers to "_FV2210"/"_FV2304"
#    (the "diff" header in between is left untouched).
# 2. Turn the A1:U63 data range into a native Excel Table ("Table1").
# 3. Freeze the header row (row 1) in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames ------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2. Convert the data range into an Excel Table ---------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the top row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
